$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cells for cell formats already present in the sheet:
#  - B2 uses style with wrapText (style index 2)
#  - C2 uses style without wrapText (style index 1)
$wrapDonor = $ws.Range("B2")
$noWrapDonor = $ws.Range("C2")

# New repeating 8-row schedule pattern for column H (rows 2-25),
# mirrors the pattern already used in column G but shifted, and
# introduces a new "Insects" entry.
$pattern = @(
  @{ Value = "Build / Build_Speed";          Wrap = $true  },
  @{ Value = "Insects";                      Wrap = $true  },
  @{ Value = "Build / Evolve / Hatch_Speed"; Wrap = $true  },
  @{ Value = "Insects";                      Wrap = $true  },
  @{ Value = "Speed up";                     Wrap = $false },
  @{ Value = "Insects";                      Wrap = $true  },
  @{ Value = "Build / Evolve / Hatch_Speed"; Wrap = $true  },
  @{ Value = "Insects";                      Wrap = $true  }
)

for ($row = 2; $row -le 25; $row++) {
    $idx = ($row - 2) % 8
    $item = $pattern[$idx]
    $cell = $ws.Range("H$row")

    if ($item.Wrap) {
        $wrapDonor.Copy()
    } else {
        $noWrapDonor.Copy()
    }
    $cell.PasteSpecial(-4122) | Out-Null

    $cell.Value = $item.Value
}

# Update the active selection to reflect the saved view state.
$ws.Range("J23").Select()
